# Auto-generated edit script: inserts three new paragraphs after the
# "...PLI has had no significant effect..." paragraph, matching the
# target diff (OLS regression block + surrounding body text).

$d = $word.ActiveDocument
$lb = [char]11   # Word "manual line break" -> serializes as <w:br/>

# --- Locate anchor paragraph and its index ---
$anchorText = "*PLI has had no significant effect on the investment levels*"
$allParas0 = $d.Paragraphs
$anchorIndex = 0
for ($i = 1; $i -le $allParas0.Count; $i++) {
    if ($allParas0.Item($i).Range.Text -like $anchorText) {
        $anchorIndex = $i
        break
    }
}
$anchorPara = $allParas0.Item($anchorIndex)
$startPos = $anchorPara.Range.End

# --- Insert 3 empty paragraph marks in a row right after the anchor ---
$rng = $anchorPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.InsertParagraphAfter()
$rng.InsertParagraphAfter()

# --- Assign paragraph styles by (stable) paragraph index ---
$allParas1 = $d.Paragraphs
$allParas1.Item($anchorIndex + 1).Style = "BodyText"
$allParas1.Item($anchorIndex + 2).Style = "SourceCode"
$allParas1.Item($anchorIndex + 3).Style = "FirstParagraph"

# --- Fill in text for each new paragraph (absolute offsets, precomputed) ---
$pos = $startPos
# paragraph 1: style=BodyText
$t0 = 'In model 3,4 and 5, the AR(1) statistic turns out to be insignificant at the 5% level, that is we fail to reject the null hypothesis that there is no autocorrelation of order 1. This means that in these models, the application of OLS should give reasonable estimates. This is what is done next.'
$r0 = $d.Range($pos, $pos)
$r0.Text = $t0
$pos = $pos + $t0.Length + 1   # +1 to skip over this paragraph's mark

# paragraph 2: style=SourceCode
$pieces1 = @(
  '## ',
  '## OLS Panel Regression Models of effect of PLI scheme on Category 1 beneficiary investment',
  '## =======================================================',
  '##                                Dependent variable:     ',
  '##                           -----------------------------',
  '##                                      i.by.k            ',
  '##                              (3)        (4)      (5)   ',
  '## -------------------------------------------------------',
  '## i_1.by.k_1                  -0.045    -0.054    -0.045 ',
  '##                            (0.053)    (0.036)  (0.052) ',
  '##                                                        ',
  '## cf_1.by.k_1                 0.004     -0.011    0.002  ',
  '##                            (0.050)    (0.032)  (0.050) ',
  '##                                                        ',
  '## ds.by.k                     0.005               0.005  ',
  '##                            (0.004)             (0.004) ',
  '##                                                        ',
  '## s_1.by.k_1                           0.017***          ',
  '##                                       (0.006)          ',
  '##                                                        ',
  '## d_1.by.k_1                  0.124*     0.090    0.127* ',
  '##                            (0.066)    (0.092)  (0.065) ',
  '##                                                        ',
  '## uncertainty_1              -0.001*    -0.001    -0.001 ',
  '##                            (0.001)    (0.001)  (0.001) ',
  '##                                                        ',
  '## repo_rate                   0.008                      ',
  '##                            (0.010)                     ',
  '##                                                        ',
  '## d.repo_rate                           -0.004    -0.003 ',
  '##                                       (0.018)  (0.019) ',
  '##                                                        ',
  '## pli                        -0.043**   -0.016    -0.038 ',
  '##                            (0.017)    (0.034)  (0.036) ',
  '##                                                        ',
  '## -------------------------------------------------------',
  '## Observations                 340        340      340   ',
  '## R2                          0.060      0.171    0.058  ',
  '## Adjusted R2                 -0.110     0.021    -0.113 ',
  '## F Statistic (df = 7; 287)  2.627**   8.448***  2.528** ',
  '## =======================================================',
  '## Note:                       *p<0.1; **p<0.05; ***p<0.01'
)
$fullText1 = [string]::Join($lb, $pieces1)
$r1 = $d.Range($pos, $pos)
$r1.Text = $fullText1
$cursor1 = $pos
foreach ($piece in $pieces1) {
    $pieceStart = $cursor1
    $pieceEnd = $pieceStart + $piece.Length
    $pieceRange = $d.Range($pieceStart, $pieceEnd)
    $pieceRange.Style = "VerbatimChar"
    $cursor1 = $pieceEnd + 1
}
$pos = $pos + $fullText1.Length + 1   # +1 to skip over this paragraph's mark

# paragraph 3: style=FirstParagraph
$t2 = 'The above results show that the model has very little explanatory power and PLI is not significant.'
$r2 = $d.Range($pos, $pos)
$r2.Text = $t2
$pos = $pos + $t2.Length + 1   # +1 to skip over this paragraph's mark

Write-Output "Inserted OLS regression section after anchor paragraph (index $anchorIndex)."